$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells: I1 = "I0", J1 = "IF" (same header style as the
# existing H1 header cell, i.e. bold / bordered / centered-top).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data rows 2..41: column I is a constant 1, column J duplicates
# whatever is already in column H for that row.
for ($r = 2; $r -le 41; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
